$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay plain text (they mirror the source site's raw
# formatting -- thousand-separator dots, fixed decimal places, etc.) -- so each
# updated cell is switched to Text format before assigning, then restored to the
# workbook's default (unstyled) look afterwards.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.288.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.427.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '659.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.445'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.08'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.423.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.990.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.064.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.607'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.438.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '519.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000209'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.610.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.154'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '30.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '540.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.155'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.889'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.04'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0441'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.72'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.37'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.13'
$ws.Range('D51').Style = 'Normal'

# Coin name / link / volume(1h) text cells
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('E3').Value = '  +2.56%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('E7').Value = '  -5.69%  '
$ws.Range('E8').Value = '  +3.18%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('E12').Value = '  +4.01%  '
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('E14').Value = '  +16.10%  '
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('E18').Value = '  +29.30%  '
$ws.Range('E19').Value = '  +36.20%  '
$ws.Range('E20').Value = '  +2.87%  '
$ws.Range('E21').Value = '  +8.20%  '
$ws.Range('E22').Value = '  +8.10%  '
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  -3.63%  '
$ws.Range('E25').Value = '  +2.32%  '
$ws.Range('E26').Value = '  +4.60%  '
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('E28').Value = '  +5.69%  '
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('E31').Value = '  +9.02%  '
$ws.Range('E32').Value = '  +5.18%  '
$ws.Range('E34').Value = '  +12.18%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +15.23%  '
$ws.Range('E37').Value = '  +3.75%  '
$ws.Range('E38').Value = '  +4.70%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E39').Value = '  +11.88%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  +19.86%  '
$ws.Range('E44').Value = '  +9.40%  '
$ws.Range('E45').Value = '  +21.95%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E47').Value = '  +9.86%  '
$ws.Range('E48').Value = '  -4.51%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E49').Value = '  +15.60%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E50').Value = '  +4.72%  '
$ws.Range('E51').Value = '  +4.71%  '
